$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 456
$ws.Range("B3").Value = "m"
$ws.Range("C3").Value = "vsdfds"
$ws.Range("D3").Value = 50
$ws.Range("E3").Value = 12
$ws.Range("F3").Value = 12
$ws.Range("G3").Value = 12

$ws.Range("G4").Select()
